# Update question identifiers from 130120 to 220120, and update the
# privateGroup label for the first question row from "Public" to
# "All Inside Track Members" (inclusion of member persona's regression).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 8; $row++) {
    $cellA = $ws.Cells.Item($row, 1)
    $oldA = $cellA.Value()
    if ($oldA -ne $null) {
        $cellA.Value = $oldA -replace '130120', '220120'
    }

    $cellB = $ws.Cells.Item($row, 2)
    $oldB = $cellB.Value()
    if ($oldB -ne $null) {
        $cellB.Value = $oldB -replace '130120', '220120'
    }
}

$ws.Range("C2").Value = "All Inside Track Members"
